# ---------------------------------------------------------------------------
# "Add download photos in Excel"
#   - Rename Sheet1 -> 報名表資料
#   - Add a new sheet 報名表照片 with a team-photo roster table
#   - Sheet1: add a "序號" header label above the numbering column, move the
#     saved selection to D16
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename the existing sheet, add the new one right after it ---------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "報名表資料"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "報名表照片"

# --- 2. 報名表照片 sheet layout --------------------------------------------

$ws2.Rows.Item(1).RowHeight = 17
$ws2.Rows.Item(2).RowHeight = 17
$ws2.Rows.Item(3).RowHeight = 17
$ws2.Rows.Item(4).RowHeight = 23
$ws2.Rows.Item(5).RowHeight = 17
$ws2.Rows.Item(6).RowHeight = 91
$ws2.Rows.Item(7).RowHeight = 15
$ws2.Rows.Item(8).RowHeight = 91
$ws2.Rows.Item(9).RowHeight = 15
$ws2.Rows.Item(10).RowHeight = 91
$ws2.Rows.Item(11).RowHeight = 15
$ws2.Rows.Item(12).RowHeight = 91
$ws2.Rows.Item(13).RowHeight = 15
$ws2.Rows.Item(14).RowHeight = 91
$ws2.Rows.Item(15).RowHeight = 15

$ws2.Range("A1:E1").ColumnWidth = 18.83

# Title
$title = $ws2.Range("C1")
$title.Font.Bold = $true
$title.Font.Size = 14
$title.Font.Name = "BiauKai"
$title.Value = "附件二：報名表照片"

# Row 2 : 報名單位 / 聯絡人
$c = $ws2.Range("A2")
$c.HorizontalAlignment = -4152
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14
$c.Value = "報名單位："

$c = $ws2.Range("B2")
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("C2")
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("D2")
$c.HorizontalAlignment = -4152
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14
$c.Value = "聯絡人："

$c = $ws2.Range("E2")
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("F2")
$c.HorizontalAlignment = -4131
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("G2")
$c.HorizontalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

# Row 3 : 電話 / 備用手機
$c = $ws2.Range("A3")
$c.HorizontalAlignment = -4152
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14
$c.Value = "電話："

$c = $ws2.Range("B3")
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("C3")
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("D3")
$c.HorizontalAlignment = -4152
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14
$c.Value = "備用手機："

$c = $ws2.Range("E3")
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("F3")
$c.HorizontalAlignment = -4131
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("G3")
$c.HorizontalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

# Row 4 : 電子郵件
$c = $ws2.Range("A4")
$c.HorizontalAlignment = -4152
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14
$c.Value = "電子郵件："

$c = $ws2.Range("B4")
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("C4")
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("D4")
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("E4")
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("F4")
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("G4")
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

# Row 5 : spacer row
$c = $ws2.Range("A5")
$c.HorizontalAlignment = -4152
$c.VerticalAlignment = -4108
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("C5")
$c.HorizontalAlignment = -4131
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("D5")
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("E5")
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("F5")
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

$c = $ws2.Range("G5")
$c.Font.Name = "KaiTi"
$c.Font.Size = 14

# Row 6 : photo placeholder box #1
$ws2.Range("A6:E6").Borders.LineStyle = 1

# Row 7 : name labels 1-5
$lbl = $ws2.Range("A7:E7")
$lbl.HorizontalAlignment = -4131
$lbl.Font.Name = "KaiTi"
$lbl.Font.Size = 12
$lbl.Borders.LineStyle = 1
$ws2.Range("A7").Value = "1.姓名："
$ws2.Range("B7").Value = "2.姓名："
$ws2.Range("C7").Value = "3.姓名："
$ws2.Range("D7").Value = "4.姓名："
$ws2.Range("E7").Value = "5.姓名："

# Row 8 : photo placeholder box #2
$ws2.Range("A8:E8").Borders.LineStyle = 1

# Row 9 : name labels 6-10
$lbl = $ws2.Range("A9:E9")
$lbl.HorizontalAlignment = -4131
$lbl.Font.Name = "KaiTi"
$lbl.Font.Size = 12
$lbl.Borders.LineStyle = 1
$ws2.Range("A9").Value = "6.姓名："
$ws2.Range("B9").Value = "7.姓名："
$ws2.Range("C9").Value = "8.姓名："
$ws2.Range("D9").Value = "9.姓名："
$ws2.Range("E9").Value = "10.姓名："

# Row 10 : photo placeholder box #3
$ws2.Range("A10:E10").Borders.LineStyle = 1

# Row 11 : name labels 11-15
$lbl = $ws2.Range("A11:E11")
$lbl.HorizontalAlignment = -4131
$lbl.Font.Name = "KaiTi"
$lbl.Font.Size = 12
$lbl.Borders.LineStyle = 1
$ws2.Range("A11").Value = "11.姓名："
$ws2.Range("B11").Value = "12.姓名："
$ws2.Range("C11").Value = "13.姓名："
$ws2.Range("D11").Value = "14.姓名："
$ws2.Range("E11").Value = "15.姓名："

# Row 12 : photo placeholder box #4
$ws2.Range("A12:E12").Borders.LineStyle = 1

# Row 13 : name labels 16-20
$lbl = $ws2.Range("A13:E13")
$lbl.HorizontalAlignment = -4131
$lbl.Font.Name = "KaiTi"
$lbl.Font.Size = 12
$lbl.Borders.LineStyle = 1
$ws2.Range("A13").Value = "16.姓名："
$ws2.Range("B13").Value = "17.姓名："
$ws2.Range("C13").Value = "18.姓名："
$ws2.Range("D13").Value = "19.姓名："
$ws2.Range("E13").Value = "20.姓名："

# Row 14 : photo placeholder box #5
$ws2.Range("A14:E14").Borders.LineStyle = 1

# Row 15 : name labels 21-25
$lbl = $ws2.Range("A15:E15")
$lbl.HorizontalAlignment = -4131
$lbl.Font.Name = "KaiTi"
$lbl.Font.Size = 12
$lbl.Borders.LineStyle = 1
$ws2.Range("A15").Value = "21.姓名："
$ws2.Range("B15").Value = "22.姓名："
$ws2.Range("C15").Value = "23.姓名："
$ws2.Range("D15").Value = "24.姓名："
$ws2.Range("E15").Value = "25.姓名："

# --- 3. 報名表資料 (Sheet1) edits ------------------------------------------

$hdr = $ws1.Range("A7")
$hdr.Font.Bold = $true
$hdr.Font.Name = "KaiTi"
$hdr.Value = "序號"

# --- 4. Selections / active sheet ------------------------------------------

$ws2.Range("C10").Select()
$ws1.Select()
$ws1.Range("D16").Select()
